# Update 'want to go' counts (column F) across sheets, per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 287
$ws.Cells.Item(3, 6).Value = 90
$ws.Cells.Item(4, 6).Value = 1204
$ws.Cells.Item(5, 6).Value = 834
$ws.Cells.Item(6, 6).Value = 865
$ws.Cells.Item(7, 6).Value = 1574
$ws.Cells.Item(8, 6).Value = 323
$ws.Cells.Item(9, 6).Value = 1075
$ws.Cells.Item(10, 6).Value = 36
$ws.Cells.Item(11, 6).Value = 83
$ws.Cells.Item(12, 6).Value = 210
$ws.Cells.Item(14, 6).Value = 538
$ws.Cells.Item(15, 6).Value = 81
$ws.Cells.Item(16, 6).Value = 50
$ws.Cells.Item(19, 6).Value = 303
$ws.Cells.Item(20, 6).Value = 594
$ws.Cells.Item(21, 6).Value = 591
$ws.Cells.Item(22, 6).Value = 71
$ws.Cells.Item(23, 6).Value = 14
$ws.Cells.Item(24, 6).Value = 793
$ws.Cells.Item(25, 6).Value = 268
$ws.Cells.Item(28, 6).Value = 382

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 18
$ws.Cells.Item(4, 6).Value = 289
$ws.Cells.Item(6, 6).Value = 201
$ws.Cells.Item(7, 6).Value = 73
$ws.Cells.Item(8, 6).Value = 601
$ws.Cells.Item(11, 6).Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 271

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 271
$ws.Cells.Item(3, 6).Value = 287
$ws.Cells.Item(5, 6).Value = 90
$ws.Cells.Item(6, 6).Value = 1204
$ws.Cells.Item(7, 6).Value = 835
$ws.Cells.Item(8, 6).Value = 865
$ws.Cells.Item(9, 6).Value = 1574
$ws.Cells.Item(10, 6).Value = 323
$ws.Cells.Item(11, 6).Value = 1075
$ws.Cells.Item(12, 6).Value = 36
$ws.Cells.Item(13, 6).Value = 83
$ws.Cells.Item(14, 6).Value = 210
$ws.Cells.Item(16, 6).Value = 538
$ws.Cells.Item(17, 6).Value = 81
$ws.Cells.Item(18, 6).Value = 50
$ws.Cells.Item(19, 6).Value = 18
$ws.Cells.Item(21, 6).Value = 289
$ws.Cells.Item(23, 6).Value = 303
$ws.Cells.Item(25, 6).Value = 201
$ws.Cells.Item(26, 6).Value = 201
$ws.Cells.Item(27, 6).Value = 594
$ws.Cells.Item(28, 6).Value = 591
$ws.Cells.Item(29, 6).Value = 71
$ws.Cells.Item(30, 6).Value = 14
$ws.Cells.Item(31, 6).Value = 793
$ws.Cells.Item(32, 6).Value = 268
$ws.Cells.Item(33, 6).Value = 73
$ws.Cells.Item(35, 6).Value = 601
$ws.Cells.Item(40, 6).Value = 382
$ws.Cells.Item(41, 6).Value = 8
